$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value = -0.8786
$ws.Range("B4").Value = -0.0348
$ws.Range("B5").Value = 0.1918
$ws.Range("B6").Value = -0.4255
$ws.Range("B7").Value = -0.6541
$ws.Range("B8").Value = 0.0443
$ws.Range("B9").Value = -0.2141
$ws.Range("B10").Value = 0.0057
$ws.Range("B11").Value = 0.0475
$ws.Range("B12").Value = -1.1992
$ws.Range("B13").Value = -0.0115
$ws.Range("B14").Value = -1.084
$ws.Range("B15").Value = -0.227
$ws.Range("B16").Value = -0.1722
$ws.Range("B17").Value = 0.0041
$ws.Range("B18").Value = 0.0404
$ws.Range("B19").Value = 0.01
$ws.Range("B20").Value = -0.2572
$ws.Range("B21").Value = 0.0175
$ws.Range("B22").Value = -0.0021
$ws.Range("B23").Value = 0.0852
$ws.Range("B24").Value = 0.0437
$ws.Range("B25").Value = -0.0092
